# Updating results with correct MI SEs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = theta_se (standard errors for theta)
$ws.Range("B4").Value = "(0.47)"
$ws.Range("F4").Value = "(0.24)"

# Row 6 = lambda_se (standard errors for lambda)
$ws.Range("B6").Value = "(0.34)"
$ws.Range("D6").Value = "(0.27)"
$ws.Range("E6").Value = "(0.28)"
$ws.Range("F6").Value = "(0.16)"
